# DailyWorkReport.xlsx - add a new daily entry block (2025-01-16) covering
# "DBFinal Demo and Database with C# basic parts complete" work, following the
# exact same 7-row layout/formatting used by the existing blocks (header row +
# "Meeting" row + continuation row + "Study" row + 3 continuation rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing block in rows 3:9 (2025-01-08) has the identical shape to the
# new block we need (header, Meeting+topic, continuation, Study+topic, then
# three more continuation/topic rows). Copying it down preserves every style
# (borders, number formats, fonts) exactly, matching how this sheet was
# originally authored (copy-paste of the previous day's block).
$src = $ws.Range("A3:D9")
$dst = $ws.Range("A36:D42")
$src.Copy($dst)

# Row 36: new day header - 2025-01-16, "Domm" entry of 0.25 hr
$ws.Range("A36").Value = 45673
$ws.Range("B36").Value = "Domm"
$ws.Range("C36").Value = ""
$ws.Range("D36").Value = 0.25

# Row 37: Meeting - General Discussion, 0.5 hr
$ws.Range("B37").Value = "Meeting"
$ws.Range("C37").Value = "General Discussion"
$ws.Range("D37").Value = 0.5

# Row 38: continuation - Reconsile, 1 hr
$ws.Range("C38").Value = "Reconsile"
$ws.Range("D38").Value = 1

# Row 39: Study - DataBase Overview query and other updates, 0.5 hr
$ws.Range("B39").Value = "Study"
$ws.Range("C39").Value = "DataBase Overview query and other updates"
$ws.Range("D39").Value = 0.5

# Row 40: continuation - Web Development Revision for Reconsile, 1.5 hr
$ws.Range("C40").Value = "Web Development Revision for Reconsile"
$ws.Range("D40").Value = 1.5

# Row 41: continuation - Database with C# (CRUD), 4 hr
$ws.Range("C41").Value = "Database with C# (CRUD)"
$ws.Range("D41").Value = 4

# Row 42: continuation - Post-meeting Discussion, 0.25 hr
$ws.Range("C42").Value = "Post-meeting Discussion"
$ws.Range("D42").Value = 0.25

# Leave the selection on the last entered cell, as in the authored workbook.
$ws.Range("D40").Select()
